$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (Exhibitions) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1339
$ws1.Range("F3").Value = 1210
$ws1.Range("F5").Value = 113
$ws1.Range("F7").Value = 672
$ws1.Range("F8").Value = 112
$ws1.Range("F11").Value = 2417
$ws1.Range("F13").Value = 1449
$ws1.Range("F14").Value = 310
$ws1.Range("F15").Value = 243
$ws1.Range("F16").Value = 590
$ws1.Range("F17").Value = 776
$ws1.Range("F18").Value = 62
$ws1.Range("F19").Value = 304
$ws1.Range("F22").Value = 24
$ws1.Range("F24").Value = 4888
$ws1.Range("F26").Value = 423
$ws1.Range("F27").Value = 70
$ws1.Range("F28").Value = 157
$ws1.Range("F29").Value = 137
$ws1.Range("F33").Value = 1034
$ws1.Range("F34").Value = 707
$ws1.Range("F36").Value = 48
$ws1.Range("F38").Value = 387
$ws1.Range("F39").Value = 1030
$ws1.Range("F40").Value = 129
$ws1.Range("F41").Value = 105

# Sheet 4: "全部类型" (All types) - mirrors the same updated values
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 1339
$ws4.Range("F5").Value = 1210
$ws4.Range("F9").Value = 113
$ws4.Range("F11").Value = 672
$ws4.Range("F12").Value = 112
$ws4.Range("F17").Value = 2417
$ws4.Range("F19").Value = 1449
$ws4.Range("F20").Value = 310
$ws4.Range("F21").Value = 243
$ws4.Range("F22").Value = 590
$ws4.Range("F25").Value = 62
$ws4.Range("F26").Value = 304
$ws4.Range("F28").Value = 24
$ws4.Range("F29").Value = 4888
$ws4.Range("F31").Value = 423
$ws4.Range("F32").Value = 70
$ws4.Range("F33").Value = 157
$ws4.Range("F34").Value = 137
$ws4.Range("F38").Value = 1034
$ws4.Range("F39").Value = 707
$ws4.Range("F40").Value = 48
$ws4.Range("F41").Value = 387
$ws4.Range("F42").Value = 1030
$ws4.Range("F43").Value = 129
$ws4.Range("F44").Value = 105
